# Fruta / hortaliza, semanal
# Insert a new weekly data row at row 126 (pushing existing rows 126-130 down to 127-131)
# and populate it with the new week's values for "Berenjena" at "Macroferia Regional de Talca".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 126, shifting rows 126:130 to 127:131.
$ws.Rows.Item(126).Insert()

# Populate the new row 126 with the new record.
$ws.Cells.Item(126, 1).Value = 5
$ws.Cells.Item(126, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(126, 3).Value = "Maule"
$ws.Cells.Item(126, 4).Value = 44747
$ws.Cells.Item(126, 5).Value = 7
$ws.Cells.Item(126, 6).Value = 100112001
$ws.Cells.Item(126, 7).Value = "Berenjena"
$ws.Cells.Item(126, 8).Value = "Sin especificar"
$ws.Cells.Item(126, 9).Value = "Primera"
$ws.Cells.Item(126, 10).Value = 300
$ws.Cells.Item(126, 11).Value = 7000
$ws.Cells.Item(126, 12).Value = 7000
$ws.Cells.Item(126, 13).Value = 7000
$ws.Cells.Item(126, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(126, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(126, 16).Value = 140
$ws.Cells.Item(126, 17).Value = 50
$ws.Cells.Item(126, 18).Value = "Hortaliza"

# Ensure the date cell keeps the date-time number format used by the other rows in column D.
$ws.Cells.Item(126, 4).NumberFormat = $ws.Cells.Item(127, 4).NumberFormat
